$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 4 (the "Search Store" test-plan row) entirely; rows below shift up.
$ws.Rows(4).Delete()

# Update the selection/view to match the post-edit state.
$ws.Range("A4:M4").Select()
